$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" column header (H1), copying the formatting used by the other
# header cells (e.g. G1: bold font, border, centered alignment)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("H1").Value = "Save"

# New "Save" column data values
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 1
